# Add data for 2022-04-05
# - Rename the "Through 2022-03-27" sheet/header to "Through 2022-03-28"
# - Update the carjacking counts for the current month (column B, "March 2022")
#   and a handful of prior-year cells picked up by the refreshed query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update its header label to reflect the new "through" date.
$ws.Name = "Through 2022-03-28"
$ws.Range("B1").Value = "March 2022 (through March 28)"

# Updated / new counts.
$ws.Range("B3").Value = 12
$ws.Range("K4").Value = 3
$ws.Range("N7").Value = 1
$ws.Range("B8").Value = 1
$ws.Range("B10").Value = 2
$ws.Range("Q10").Value = 2
$ws.Range("N11").Value = 7
$ws.Range("Q11").Value = 1
$ws.Range("B12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("W15").Value = 4
$ws.Range("B24").Value = 3
$ws.Range("B27").Value = 3
$ws.Range("B32").Value = 4
$ws.Range("B55").Value = 2
$ws.Range("B62").Value = 1
$ws.Range("W90").Value = 1
